# Fixed a bug in countlines
# The rows of data (A2:F21) need to be reordered. The commit fixes a bug
# in the "countlines" logic, which results in the rows being written out
# in a different (corrected) order, while the totals row (26) and header
# row (1) stay the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for rows 2..21, columns A..F
$data = @{
    2  = @(901, 16, 15, 45, 60, 60)
    3  = @(1202, 2, 10, 10, 10, 10)
    4  = @(501, 9, 52, 30, 75, 45)
    5  = @(601, 9, 60, 67, 60, 42)
    6  = @(201, 9, 30, 15, 45, 30)
    7  = @(801, 3, 67, 65, 52, 45)
    8  = @(1201, 2, 10, 10, 10, 10)
    9  = @(101, 9, 30, 15, 60, 15)
    10 = @(1001, 18, 30, 75, 60, 72)
    11 = @(401, 9, 48, 67, 75, 45)
    12 = @(701, 3, 90, 45, 97, 15)
    13 = @(1203, 3, 15, 15, 15, 15)
    14 = @(902, 1, 0, 0, 0, 0)
    15 = @(301, 6, 45, 30, 60, 45)
    16 = @(3, 0, 3, 3, 3, 3)
    17 = @(1101, 0, 15, 30, 30, 0)
    18 = @(802, 0, 4, 5, 4, 0)
    19 = @(1, 0, 2, 2, 2, 2)
    20 = @(502, 0, 4, 0, 0, 0)
    21 = @(2, 0, 2, 2, 2, 2)
}

$cols = @("A", "B", "C", "D", "E", "F")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt 6; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $values[$i]
    }
}
